$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 33 with a new mac-address/document-type test record
$ws.Range("A33").Value = 10002
$ws.Range("B33").Value = 110032
$ws.Range("C33").Value = 10032
$ws.Range("D33").Value = "eng"
$ws.Range("E33").Value = $true
$ws.Range("F33").Value = "superadmin"
$ws.Range("G33").Value = "now()"

# Update selection to reflect the new active cell after entry
$ws.Range("C30").Select()
